$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 56.625
$ws.Range("I11").Value = 56.625
$ws.Range("K11").Value = 56.625
$ws.Range("M11").Value = 83.375
$ws.Range("H32").Value = 15376
$ws.Range("I32").Value = 5500
$ws.Range("J32").Value = 18668
$ws.Range("K32").Value = 5500
$ws.Range("L32").Value = 18668
$ws.Range("M32").Value = -5174
$ws.Range("N32").Value = -19320
$ws.Range("H33").Value = 522.9091
$ws.Range("I33").Value = 167.23529
$ws.Range("J33").Value = 1732.2
$ws.Range("K33").Value = 167.23529
$ws.Range("L33").Value = 1732.2
$ws.Range("M33").Value = 61.76471000000001
$ws.Range("N33").Value = -2190.2
$ws.Range("H74").Value = 9828
$ws.Range("I74").Value = 8821.875
$ws.Range("K74").Value = 8821.875
$ws.Range("M74").Value = -7885.875
$ws.Range("H77").Value = 9828
$ws.Range("I77").Value = 8821.875
$ws.Range("K77").Value = 44109.375
$ws.Range("M77").Value = -39429.375
$ws.Range("H80").Value = 2487.5
$ws.Range("J80").Value = 4009.5557
$ws.Range("L80").Value = 12028.6671
$ws.Range("N80").Value = -14024.6671
$ws.Range("H83").Value = 2487.5
$ws.Range("J83").Value = 4009.5557
$ws.Range("L83").Value = 36086.0013
$ws.Range("N83").Value = -46070.0013
$ws.Range("H107").Value = 1863.0714
$ws.Range("I107").Value = 2034.35
$ws.Range("J107").Value = 1434.875
$ws.Range("K107").Value = 2034.35
$ws.Range("L107").Value = 1434.875
$ws.Range("M107").Value = -114.3499999999999
$ws.Range("N107").Value = -5274.875
$ws.Range("H135").Value = 661.6957
$ws.Range("I135").Value = 581.8570999999999
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 5236.7139
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -2701.7139
$ws.Range("N135").Value = -18570
$ws.Range("H138").Value = 2523.2092
$ws.Range("J138").Value = 3208.35
$ws.Range("L138").Value = 9625.049999999999
$ws.Range("N138").Value = -19905.05
$ws.Range("H141").Value = 6166.5835
$ws.Range("I141").Value = 1333
$ws.Range("J141").Value = 7777.778
$ws.Range("K141").Value = 3999
$ws.Range("L141").Value = 23333.334
$ws.Range("M141").Value = 1181
$ws.Range("N141").Value = -33693.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 232.5
$ws.Range("I4").Value = 232.5
$ws.Range("K4").Value = 232.5
$ws.Range("M4").Value = -116.5
$ws.Range("H32").Value = 5885.1763
$ws.Range("I32").Value = 5105.6
$ws.Range("J32").Value = 11732
$ws.Range("K32").Value = 5105.6
$ws.Range("L32").Value = 11732
$ws.Range("M32").Value = -4818.6
$ws.Range("N32").Value = -12306
$ws.Range("H61").Value = 5132.0454
$ws.Range("I61").Value = 3783.7222
$ws.Range("K61").Value = 3783.7222
$ws.Range("M61").Value = -3571.7222
$ws.Range("H74").Value = 22224786
$ws.Range("I74").Value = 25643676
$ws.Range("K74").Value = 25643676
$ws.Range("M74").Value = -25642802
$ws.Range("H77").Value = 22224786
$ws.Range("I77").Value = 25643676
$ws.Range("K77").Value = 128218380
$ws.Range("M77").Value = -128214012
$ws.Range("H102").Value = 1919
$ws.Range("I102").Value = 1996.5
$ws.Range("K102").Value = 1996.5
$ws.Range("M102").Value = -374.5
$ws.Range("H122").Value = 2080.6538
$ws.Range("I122").Value = 1610.3684
$ws.Range("K122").Value = 4831.1052
$ws.Range("M122").Value = -2381.1052
$ws.Range("H132").Value = 3272.3333
$ws.Range("I132").Value = 2449.25
$ws.Range("K132").Value = 7347.75
$ws.Range("M132").Value = -4817.75
$ws.Range("H136").Value = 5132.0454
$ws.Range("I136").Value = 3783.7222
$ws.Range("K136").Value = 11351.1666
$ws.Range("M136").Value = -8801.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1029.1364
$ws.Range("I94").Value = 942.15
$ws.Range("K94").Value = 942.15
$ws.Range("M94").Value = -491.15
$ws.Range("H107").Value = 1108.6
$ws.Range("I107").Value = 1108.6
$ws.Range("K107").Value = 1108.6
$ws.Range("M107").Value = 811.4000000000001
$ws.Range("H134").Value = 5220.5
$ws.Range("I134").Value = 1577.75
$ws.Range("K134").Value = 4733.25
$ws.Range("M134").Value = -2198.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 595.5185
$ws.Range("J7").Value = 381.27274
$ws.Range("L7").Value = 381.27274
$ws.Range("N7").Value = -607.27274
$ws.Range("H31").Value = 31579.027
$ws.Range("I31").Value = 3291.7036
$ws.Range("K31").Value = 3291.7036
$ws.Range("M31").Value = -2996.7036
$ws.Range("H34").Value = 31579.027
$ws.Range("I34").Value = 3291.7036
$ws.Range("K34").Value = 3291.7036
$ws.Range("M34").Value = -3089.7036

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 980
$ws.Range("J23").Value = 1343.7
$ws.Range("L23").Value = 4031.1
$ws.Range("N23").Value = -4501.1
$ws.Range("H98").Value = 999.25
$ws.Range("I98").Value = 998.5
$ws.Range("K98").Value = 2995.5
$ws.Range("M98").Value = -1497.5
$ws.Range("H133").Value = 1883.3334
$ws.Range("I133").Value = 1883.3334
$ws.Range("K133").Value = 5650.0002
$ws.Range("M133").Value = -590.0002000000004
$ws.Range("H138").Value = 47510
$ws.Range("J138").Value = 16250
$ws.Range("L138").Value = 48750
$ws.Range("N138").Value = -59030
$ws.Range("H139").Value = 6509.625
$ws.Range("I139").Value = 2174.1667
$ws.Range("J139").Value = 19516
$ws.Range("K139").Value = 6522.500100000001
$ws.Range("L139").Value = 58548
$ws.Range("M139").Value = -1382.500100000001
$ws.Range("N139").Value = -68828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value = 0
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("N82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("N85").Value = 0
$ws.Range("H122").Value = 3191.35
$ws.Range("I122").Value = 2582.4375
$ws.Range("K122").Value = 7747.3125
$ws.Range("M122").Value = -5297.3125
$ws.Range("H126").Value = 4035.0667
$ws.Range("I126").Value = 2917.6667
$ws.Range("K126").Value = 8753.000100000001
$ws.Range("M126").Value = -6283.000100000001
$ws.Range("H132").Value = 6912.875
$ws.Range("I132").Value = 3657.8
$ws.Range("J132").Value = 12338
$ws.Range("K132").Value = 10973.4
$ws.Range("L132").Value = 37014
$ws.Range("M132").Value = -8443.400000000001
$ws.Range("N132").Value = -42074

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7462.696
$ws.Range("I40").Value = 7006.6
$ws.Range("K40").Value = 7006.6
$ws.Range("M40").Value = -6870.6
$ws.Range("H93").Value = 2152.9412
$ws.Range("I93").Value = 1719.08
$ws.Range("K93").Value = 1719.08
$ws.Range("M93").Value = -471.0799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 12101.8
$ws.Range("I64").Value = 12101.8
$ws.Range("K64").Value = 12101.8
$ws.Range("M64").Value = -11853.8
$ws.Range("H67").Value = 12101.8
$ws.Range("I67").Value = 12101.8
$ws.Range("K67").Value = 12101.8
$ws.Range("M67").Value = -11243.8
$ws.Range("H113").Value = 1008
$ws.Range("J113").Value = 1008
$ws.Range("L113").Value = 3024
$ws.Range("N113").Value = -7364
$ws.Range("H132").Value = 4531.9062
$ws.Range("I132").Value = 4173.1724
$ws.Range("J132").Value = 7999.6665
$ws.Range("K132").Value = 12519.5172
$ws.Range("L132").Value = 23998.9995
$ws.Range("M132").Value = -9989.517200000002
$ws.Range("N132").Value = -29058.9995
